# Weekly update: insert a new week's price data (2 rows) at the top of the
# existing Limón price series for Vega Monumental Concepción, shifting the
# previously-entered rows down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 172, pushing rows 172:249
# down to 174:251 (dimension grows from A1:T249 to A1:T251).
$ws.Range("A172:T173").Insert()

# Populate the newly inserted row 172 ("1a amarillo").
$ws.Cells.Item(172, 1).Value = 11
$ws.Cells.Item(172, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(172, 3).Value = "Bíobío"
$ws.Cells.Item(172, 4).Value = 44460
$ws.Cells.Item(172, 5).Value = 8
$ws.Cells.Item(172, 6).Value = "Fruta"
$ws.Cells.Item(172, 7).Value = 100102
$ws.Cells.Item(172, 8).Value = "Cítricos"
$ws.Cells.Item(172, 9).Value = 100102003
$ws.Cells.Item(172, 10).Value = "Limón"
$ws.Cells.Item(172, 11).Value = "Sin especificar"
$ws.Cells.Item(172, 12).Value = "1a amarillo"
$ws.Cells.Item(172, 13).Value = 300
$ws.Cells.Item(172, 14).Value = 7000
$ws.Cells.Item(172, 15).Value = 7000
$ws.Cells.Item(172, 16).Value = 7000
$ws.Cells.Item(172, 17).Value = "$/malla 16 kilos"
$ws.Cells.Item(172, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(172, 19).Value = 438
$ws.Cells.Item(172, 20).Value = 16

# Populate the newly inserted row 173 ("2a amarillo").
$ws.Cells.Item(173, 1).Value = 11
$ws.Cells.Item(173, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(173, 3).Value = "Bíobío"
$ws.Cells.Item(173, 4).Value = 44460
$ws.Cells.Item(173, 5).Value = 8
$ws.Cells.Item(173, 6).Value = "Fruta"
$ws.Cells.Item(173, 7).Value = 100102
$ws.Cells.Item(173, 8).Value = "Cítricos"
$ws.Cells.Item(173, 9).Value = 100102003
$ws.Cells.Item(173, 10).Value = "Limón"
$ws.Cells.Item(173, 11).Value = "Sin especificar"
$ws.Cells.Item(173, 12).Value = "2a amarillo"
$ws.Cells.Item(173, 13).Value = 300
$ws.Cells.Item(173, 14).Value = 6000
$ws.Cells.Item(173, 15).Value = 6000
$ws.Cells.Item(173, 16).Value = 6000
$ws.Cells.Item(173, 17).Value = "$/malla 16 kilos"
$ws.Cells.Item(173, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(173, 19).Value = 375
$ws.Cells.Item(173, 20).Value = 16
